$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B1").Value = 21
$ws.Range("B2").Value = 42
$ws.Range("B3").Value = 66
$ws.Range("B4").Value = 122
$ws.Range("B5").Value = 173
$ws.Range("B6").Value = 214
$ws.Range("B7").Value = 263
$ws.Range("B8").Value = 288
$ws.Range("B9").Value = 355
$ws.Range("B10").Value = 381
